$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (previously example "Иванов / Иван / Иванович")
$ws.Range("A2").Value = "Федоров"
$ws.Range("B2").Value = "Кирилл"
$ws.Range("C2").Value = "Евгеньевич"

# Update row 3 (previously example "Тестов / Тест") and add new C3 value
$ws.Range("A3").Value = "A"
$ws.Range("B3").Value = "B"
$ws.Range("C3").Value = "C"

# Update the active selection to match the saved state
$ws.Range("C6").Select()
